# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D) for the row corresponding
# to the 17c2073f-ec69-4fe4-91f4-6ccf4988d3d9 file on both the zh-cn and
# de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-15 02:35:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-15 02:35:50"
